# Adds the "July 2nd" (2020-07-02, Excel serial 44014) raw + clean SSA data
# to the bitacora_historica_datos_abiertos workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: out_vars - append row 33 with the new day's aggregated values
# ---------------------------------------------------------------------------
$wsOut = $wb.Worksheets.Item("out_vars")

$wsOut.Range("A32:J32").Copy()
$wsOut.Range("A33:J33").PasteSpecial(-4122)

$wsOut.Range("A33").Value = 44014
$wsOut.Range("B33").Value = 238511
$wsOut.Range("C33").Value = 295561
$wsOut.Range("D33").Value = 76423
$wsOut.Range("E33").Value = 29189
$wsOut.Range("F33").Value = 30.515573705195987
$wsOut.Range("G33").Value = 72783
$wsOut.Range("H33").Value = 6214
$wsOut.Range("I33").Value = 7080
$wsOut.Range("J33").Value = 610495

$wsOut.Range("A33").Select()

# ---------------------------------------------------------------------------
# Sheet 2: dates_dx - fill in the previously-blank placeholder row 33
# ---------------------------------------------------------------------------
$wsDx = $wb.Worksheets.Item("dates_dx")

$wsDx.Range("A33").Value = 44014
$wsDx.Range("B33").Value = 0
$wsDx.Range("C33").Value = 1
$wsDx.Range("D33").Value = 1
$wsDx.Range("E33").Value = 1
$wsDx.Range("F33").Value = 1
$wsDx.Range("G33").Value = 0
$wsDx.Range("H33").Value = 1
$wsDx.Range("I33").Value = 0
$wsDx.Range("J33").Value = 1
$wsDx.Range("K33").Value = 0
$wsDx.Range("L33").Value = 4

$wsDx.Range("B35").Select()

# ---------------------------------------------------------------------------
# Sheet 3: dates_sx - finish row 32 and fill in row 33
# ---------------------------------------------------------------------------
$wsSx = $wb.Worksheets.Item("dates_sx")

$wsSx.Range("D32").Value = 1
$wsSx.Range("E32").Value = 0
$wsSx.Range("F32").Value = 1
$wsSx.Range("G32").Value = 1
$wsSx.Range("H32").Value = 1
$wsSx.Range("I32").Value = 0
$wsSx.Range("J32").Value = 1
$wsSx.Range("K32").Value = 1
$wsSx.Range("L32").Value = 1
$wsSx.Range("M32").Value = 0
$wsSx.Range("N32").Value = 0

$wsSx.Range("A33").Value = 44014
$wsSx.Range("B33").Value = 0
$wsSx.Range("C33").Value = 1
$wsSx.Range("D33").Value = 1
$wsSx.Range("E33").Value = 0
$wsSx.Range("F33").Value = 1
$wsSx.Range("G33").Value = 1
$wsSx.Range("H33").Value = 1
$wsSx.Range("I33").Value = 0
$wsSx.Range("J33").Value = 1
$wsSx.Range("K33").Value = 1
$wsSx.Range("L33").Value = 1
$wsSx.Range("M33").Value = 0
$wsSx.Range("N33").Value = 0

$wsSx.Range("A33").Select()

# ---------------------------------------------------------------------------
# Sheet 4: dates_deaths - fill in the previously-blank placeholder row 33
# ---------------------------------------------------------------------------
$wsDeaths = $wb.Worksheets.Item("dates_deaths")

$wsDeaths.Range("A33").Value = 44014
$wsDeaths.Range("B33").Value = 0
$wsDeaths.Range("C33").Value = 0
$wsDeaths.Range("D33").Value = 0
$wsDeaths.Range("E33").Value = 0
$wsDeaths.Range("F33").Value = 2
$wsDeaths.Range("G33").Value = 1
$wsDeaths.Range("H33").Value = 1
$wsDeaths.Range("I33").Value = 1
$wsDeaths.Range("J33").Value = 2

$wsDeaths.Range("D37").Select()

# ---------------------------------------------------------------------------
# Sheet 5: control_obs - new date column AG + totals, becomes the active tab
# ---------------------------------------------------------------------------
$wsControl = $wb.Worksheets.Item("control_obs")

$wsControl.Range("AG1").Value = 44014
$wsControl.Range("AG2").Value = 4041
$wsControl.Range("AG3").Value = 3856
$wsControl.Range("AG4").Value = 3856
$wsControl.Range("AG5").Value = 3856
$wsControl.Range("AG6").Value = 3856
$wsControl.Range("AG7").Value = 3047
$wsControl.Range("AG8").Value = 5678
$wsControl.Range("AG10").Value = 172
$wsControl.Range("AG11").Value = 172
$wsControl.Range("AG12").Value = 172
$wsControl.Range("AG13").Value = 172
$wsControl.Range("AG14").Value = 172
$wsControl.Range("AG15").Value = 107
$wsControl.Range("AG16").Value = 184
$wsControl.Range("AG18").Value = 933

$wsControl.Range("AE20:AE20").Copy()
$wsControl.Range("AF20:AG20").PasteSpecial(-4163)

$wsControl.Range("AH17").Select()
$wsControl.Activate()
